# Updates the cryptos list data (price + 1h volume change columns) to reflect
# the latest scrape, as well as re-ordering RenderToken/WEMIXToken (rows 46-47)
# whose ranking swapped position.
#
# Because the headless Excel engine auto-detects purely-numeric-looking text
# (e.g. "1.007") and silently stores it as a real number, we force such values
# to stay as text by prefixing the formula with a leading apostrophe - exactly
# like typing '1.007 into Excel's formula bar. Values that are not ambiguous
# (contain letters, multiple dots, spaces, %, etc.) are simply assigned as-is.

function Set-TextValue($range, $text) {
    if ($text -match '^[+-]?[0-9]+(\.[0-9]+)?$') {
        $range.Formula = "'" + $text
    } else {
        $range.Value = $text
    }
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "28.453.63"
Set-TextValue $ws.Range("E2") "  +0.27%  "

# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "1.864.27"
Set-TextValue $ws.Range("E3") "  +0.64%  "

# Row 4 - TetherUSD
Set-TextValue $ws.Range("D4") "1.007"
Set-TextValue $ws.Range("E4") "  +0.22%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "324.50"
Set-TextValue $ws.Range("E5") "  -0.34%  "

# Row 6 - USDC
Set-TextValue $ws.Range("D6") "1.005"
Set-TextValue $ws.Range("E6") "  +0.00%  "

# Row 7 - XRP
Set-TextValue $ws.Range("D7") "0.4549"
Set-TextValue $ws.Range("E7") "  -1.97%  "

# Row 8 - Cardano
Set-TextValue $ws.Range("E8") "  -1.31%  "

# Row 9 - Dogecoin
Set-TextValue $ws.Range("D9") "0.07812"
Set-TextValue $ws.Range("E9") "  -0.80%  "

# Row 10 - Polygon
Set-TextValue $ws.Range("D10") "0.9875"
Set-TextValue $ws.Range("E10") "  +2.11%  "

# Row 11 - Solana
Set-TextValue $ws.Range("D11") "21.54"
Set-TextValue $ws.Range("E11") "  -2.77%  "

# Row 12 - WrappedEther
Set-TextValue $ws.Range("D12") "1.878.86"
Set-TextValue $ws.Range("E12") "  +2.21%  "

# Row 13 - Chainlink
Set-TextValue $ws.Range("E13") "  -0.21%  "

# Row 14 - Polkadot
Set-TextValue $ws.Range("D14") "5.617"
Set-TextValue $ws.Range("E14") "  -1.61%  "

# Row 15 - TRON
Set-TextValue $ws.Range("D15") "0.06890"
Set-TextValue $ws.Range("E15") "  -0.57%  "

# Row 16 - Litecoin
Set-TextValue $ws.Range("D16") "86.72"
Set-TextValue $ws.Range("E16") "  -2.54%  "

# Row 17 - BinanceUSD
Set-TextValue $ws.Range("D17") "1.006"
Set-TextValue $ws.Range("E17") "  +0.19%  "

# Row 18 - ShibaInu
Set-TextValue $ws.Range("D18") "0.000009928"
Set-TextValue $ws.Range("E18") "  -0.64%  "

# Row 19 - Avalanche
Set-TextValue $ws.Range("D19") "16.63"
Set-TextValue $ws.Range("E19") "  -1.01%  "

# Row 20 - Dai
Set-TextValue $ws.Range("E20") "  -0.27%  "

# Row 21 - WrappedBTC
Set-TextValue $ws.Range("D21") "28.454.20"
Set-TextValue $ws.Range("E21") "  +0.38%  "

# Row 22 - Uniswap
Set-TextValue $ws.Range("D22") "5.240"
Set-TextValue $ws.Range("E22") "  -1.59%  "

# Row 23 - Cosmos
Set-TextValue $ws.Range("E23") "  -1.89%  "

# Row 24 - Toncoin
Set-TextValue $ws.Range("E24") "  -0.55%  "

# Row 25 - Wrapped liquid staked Ether 2.0
Set-TextValue $ws.Range("D25") "2.082.16"
Set-TextValue $ws.Range("E25") "  +1.55%  "

# Row 26 - Monero
Set-TextValue $ws.Range("D26") "153.69"
Set-TextValue $ws.Range("E26") "  -0.66%  "

# Row 27 - EthereumClassic
Set-TextValue $ws.Range("D27") "19.09"
Set-TextValue $ws.Range("E27") "  -0.87%  "

# Row 28 - InternetComputer (DFINITY)
Set-TextValue $ws.Range("D28") "5.665"
Set-TextValue $ws.Range("E28") "  -1.87%  "

# Row 29 - BitcoinCash
Set-TextValue $ws.Range("D29") "117.36"
Set-TextValue $ws.Range("E29") "  -1.60%  "

# Row 30 - LidoDAOToken
Set-TextValue $ws.Range("D30") "1.914"
Set-TextValue $ws.Range("E30") "  -3.39%  "

# Row 31 - Stellar
Set-TextValue $ws.Range("D31") "0.09263"
Set-TextValue $ws.Range("E31") "  -0.16%  "

# Row 32 - ImmutableX
Set-TextValue $ws.Range("D32") "0.9050"
Set-TextValue $ws.Range("E32") "  -2.98%  "

# Row 33 - Filecoin
Set-TextValue $ws.Range("D33") "5.256"
Set-TextValue $ws.Range("E33") "  -0.73%  "

# Row 34 - ARBITRUM
Set-TextValue $ws.Range("E34") "  -1.22%  "

# Row 35 - HuobiToken
Set-TextValue $ws.Range("D35") "3.293"
Set-TextValue $ws.Range("E35") "  -1.15%  "

# Row 36 - Hedera
Set-TextValue $ws.Range("D36") "0.05690"
Set-TextValue $ws.Range("E36") "  -2.41%  "

# Row 37 - TrustWalletToken
Set-TextValue $ws.Range("D37") "1.142"
Set-TextValue $ws.Range("E37") "  +0.07%  "

# Row 38 - VeChain
Set-TextValue $ws.Range("D38") "0.02052"
Set-TextValue $ws.Range("E38") "  -3.12%  "

# Row 39 - FraxShare
Set-TextValue $ws.Range("D39") "7.651"
Set-TextValue $ws.Range("E39") "  -2.00%  "

# Row 40 - TheSandbox
Set-TextValue $ws.Range("D40") "0.5548"
Set-TextValue $ws.Range("E40") "  -1.39%  "

# Row 41 - Algorand
Set-TextValue $ws.Range("D41") "0.1767"
Set-TextValue $ws.Range("E41") "  -0.14%  "

# Row 42 - Aptos
Set-TextValue $ws.Range("D42") "9.626"
Set-TextValue $ws.Range("E42") "  -3.34%  "

# Row 43 - Cronos
Set-TextValue $ws.Range("D43") "0.07100"
Set-TextValue $ws.Range("E43") "  -1.78%  "

# Row 44 - EnergySwap
Set-TextValue $ws.Range("D44") "11.55"
Set-TextValue $ws.Range("E44") "  -1.54%  "

# Row 45 - Decentraland
Set-TextValue $ws.Range("D45") "0.5228"
Set-TextValue $ws.Range("E45") "  -1.30%  "

# Row 46 - was WEMIXToken, now RenderToken (rank swap with row 47)
Set-TextValue $ws.Range("B46") "RenderToken"
Set-TextValue $ws.Range("C46") "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D46") "2.120"
Set-TextValue $ws.Range("E46") "  -1.13%  "

# Row 47 - was RenderToken, now WEMIXToken (rank swap with row 46)
Set-TextValue $ws.Range("B47") "WEMIXToken"
Set-TextValue $ws.Range("C47") "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue $ws.Range("D47") "1.128"
Set-TextValue $ws.Range("E47") "  -0.63%  "

# Row 48 - NEARProtocol
Set-TextValue $ws.Range("D48") "1.808"
Set-TextValue $ws.Range("E48") "  -1.86%  "

# Row 49 - Quant
Set-TextValue $ws.Range("D49") "111.87"
Set-TextValue $ws.Range("E49") "  -2.08%  "

# Row 50 - MXToken
Set-TextValue $ws.Range("D50") "2.426"
Set-TextValue $ws.Range("E50") "  +3.86%  "

# Row 51 - PaxDollar
Set-TextValue $ws.Range("D51") "1.005"
Set-TextValue $ws.Range("E51") "  -0.02%  "
